$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new text value (matches diff: inline strings, so force text assignment)
$updates = @{
    'D2' = '27.024.55'
    'E2' = '  -0.99%  '
    'D3' = '1.827.21'
    'E3' = '  -0.23%  '
    'D5' = '310.60'
    'E5' = '  -1.29%  '
    'E6' = '  -0.32%  '
    'D7' = '0.4628'
    'E7' = '  -2.18%  '
    'D8' = '0.3723'
    'E8' = '  +0.99%  '
    'D9' = '0.07253'
    'E9' = '  -2.52%  '
    'D10' = '0.8630'
    'E10' = '  -2.55%  '
    'D11' = '19.94'
    'E11' = '  -2.60%  '
    'D12' = '0.07810'
    'E12' = '  +6.42%  '
    'D13' = '1.848.59'
    'E13' = '  -3.04%  '
    'D14' = '5.358'
    'E14' = '  -1.31%  '
    'D15' = '6.537'
    'E15' = '  -0.35%  '
    'E16' = '  -2.51%  '
    'E17' = '  -0.11%  '
    'D18' = '0.000008700'
    'E18' = '  -1.04%  '
    'E19' = '  -0.29%  '
    'D20' = '27.139.36'
    'E20' = '  -1.82%  '
    'D21' = '14.55'
    'E21' = '  -1.49%  '
    'D22' = '5.158'
    'E22' = '  -2.42%  '
    'E23' = '  -1.03%  '
    'D24' = '2.077.67'
    'E24' = '  -1.60%  '
    'D25' = '153.13'
    'E25' = '  +0.82%  '
    'D26' = '1.837'
    'E26' = '  -2.93%  '
    'D27' = '18.22'
    'E27' = '  -2.27%  '
    'D28' = '2.094'
    'E28' = '  -2.26%  '
    'D29' = '5.141'
    'E29' = '  -1.68%  '
    'D30' = '115.13'
    'E30' = '  -1.76%  '
    'E31' = '  -1.68%  '
    'D32' = '2.963'
    'E32' = '  +0.47%  '
    'D33' = '4.445'
    'E33' = '  -2.07%  '
    'D34' = '0.7262'
    'E34' = '  -3.33%  '
    'D35' = '1.134'
    'E35' = '  -3.50%  '
    'B36' = 'RenderToken'
    'C36' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D36' = '2.465'
    'E36' = '  +1.44%  '
    'B37' = 'TrustWalletToken'
    'C37' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D37' = '1.077'
    'E37' = '  -1.28%  '
    'E38' = '  -1.97%  '
    'D39' = '0.01940'
    'D40' = '2.948'
    'E40' = '  -0.63%  '
    'D41' = '7.228'
    'E41' = '  -0.29%  '
    'D42' = '0.5166'
    'E42' = '  -2.34%  '
    'D43' = '0.1628'
    'E43' = '  -1.82%  '
    'D44' = '0.8582'
    'E44' = '  -15.15%  '
    'D45' = '8.189'
    'E45' = '  -3.54%  '
    'D46' = '0.4815'
    'E46' = '  -2.34%  '
    'D47' = '1.007'
    'E47' = '  -0.38%  '
    'D48' = '10.18'
    'E48' = '  -3.51%  '
    'D49' = '102.79'
    'E49' = '  -2.26%  '
    'D50' = '0.06259'
    'E50' = '  -0.65%  '
    'D51' = '1.622'
    'E51' = '  -2.76%  '
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
}
